# Applies the "Översikt ÅSELE" update:
#  - The "Förändrad" (C) date for every data row moves from 2023-09-10 (45179)
#    to 2023-09-11 (45180).
#  - Row 2 (previously "A 32529-2022" / SCA) and row 3 (previously
#    "A 64788-2019" / Kyrkan) swap places, and the Kyrkan row (now row 2)
#    picks up three newly-found species (Gammelgransskål, Granticka,
#    Mjölsvärting), which bumps its NT / Rödlistade / Alla arter counters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 454
$newDate = 45180

# 1) Bump the "Förändrad" date column (C) for every data row (2..454).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}

# The workbook's save/round-trip otherwise turns untouched, style-only empty
# "Artnamn" cells (rows with no species listed) into empty-string cells
# instead of leaving them blank; explicitly re-blank them so they stay
# untouched, matching the source.
$ws.Range("R49:R454").Value = $null

# 2) Move the Kyrkan case ("A 64788-2019") into row 2 with its updated
#    species/count data.
$ws.Range("A2").Value = "A 64788-2019"
$ws.Range("B2").Value = 43801
$ws.Range("D2").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E2").Value = "ÅSELE"
$ws.Range("F2").Value = "Kyrkan"
$ws.Range("G2").Value = 20.6
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 16
$ws.Range("R2").Value = "Blå taggsvamp`r`nDvärgbägarlav`r`nGammelgransskål`r`nGarnlav`r`nGranticka`r`nMjölsvärting`r`nMörk kolflarnlav`r`nRosenticka`r`nSpillkråka`r`nTretåig hackspett`r`nUllticka`r`nVaddporing`r`nVedskivlav`r`nDropptaggsvamp`r`nLuddlav`r`nTrådticka"

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/artfynd/A 64788-2019.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/kartor/A 64788-2019.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/klagomål/A 64788-2019.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/klagomålsmail/A 64788-2019.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/tillsyn/A 64788-2019.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/tillsynsmail/A 64788-2019.docx")'

# 3) Move the SCA case ("A 32529-2022") into row 3, unchanged apart from
#    its new position and the refreshed date handled above.
$ws.Range("A3").Value = "A 32529-2022"
$ws.Range("B3").Value = 44782
$ws.Range("D3").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E3").Value = "ÅSELE"
$ws.Range("F3").Value = "SCA"
$ws.Range("G3").Value = 13.6
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 12
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 14
$ws.Range("R3").Value = "Storspov`r`nTornseglare`r`nGräddporing`r`nGammelgransskål`r`nGranticka`r`nKolflarnlav`r`nLunglav`r`nSkrovellav`r`nSpillkråka`r`nTalltita`r`nTretåig hackspett`r`nUllticka`r`nLuddlav`r`nStuplav"

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/artfynd/A 32529-2022.xlsx")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/kartor/A 32529-2022.png")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/klagomål/A 32529-2022.docx")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/klagomålsmail/A 32529-2022.docx")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/tillsyn/A 32529-2022.docx")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ASELE/tillsynsmail/A 32529-2022.docx")'

# Re-pasting the wrapped, multi-line species text makes Excel auto-fit the
# row height; restore the original fixed 15pt custom height for both rows.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
